$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("C2").Value = 168
$ws.Range("E2").Value = 72.76500000000001
$ws.Range("C3").Value = 373
$ws.Range("E3").Value = 156.20121
$ws.Range("C4").Value = 164
$ws.Range("E4").Value = 124.54785
$ws.Range("C5").Value = 142
$ws.Range("E5").Value = 121.178322
$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 91.36799999999999
$ws.Range("C7").Value = 0
$ws.Range("E7").Value = 199.55025
$ws.Range("C8").Value = 0
$ws.Range("E8").Value = 150.49071
$ws.Range("C9").Value = 0
$ws.Range("E9").Value = 228.67596
$ws.Range("C10").Value = 0
$ws.Range("E10").Value = 41.877
$ws.Range("C11").Value = 0
$ws.Range("E11").Value = 73.791
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = 49.19399999999999
$ws.Range("C13").Value = 0
$ws.Range("E13").Value = 93.40649999999999
$ws.Range("C14").Value = 0
$ws.Range("E14").Value = 144.7875
$ws.Range("C15").Value = 0
$ws.Range("E15").Value = 91.36799999999999
$ws.Range("C16").Value = 0
$ws.Range("E16").Value = 57.21921
$ws.Range("C17").Value = 0
$ws.Range("E17").Value = 87.79940999999999
$ws.Range("C18").Value = 0
$ws.Range("E18").Value = 199.55025

$ws = $wb.Worksheets.Item("ForecastedAverageDemand")
$ws.Range("B2").Value = 72
$ws.Range("C2").Value = 141
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 111
$ws.Range("B3").Value = 75
$ws.Range("C3").Value = 215
$ws.Range("D3").Value = 69
$ws.Range("E3").Value = 26
$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 84
$ws.Range("E4").Value = 0
$ws.Range("B5").Value = 83
$ws.Range("D5").Value = 89
$ws.Range("E5").Value = 38
$ws.Range("B6").Value = 52
$ws.Range("C6").Value = 35
$ws.Range("D6").Value = 82
$ws.Range("E6").Value = 80
$ws.Range("B7").Value = 74
$ws.Range("C7").Value = 45
$ws.Range("D7").Value = 20
$ws.Range("E7").Value = 189
$ws.Range("B8").Value = 46
$ws.Range("C8").Value = 111
$ws.Range("D8").Value = 41
$ws.Range("E8").Value = 47
$ws.Range("B9").Value = 48
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 125
$ws.Range("E9").Value = 22
$ws.Range("B10").Value = 84
$ws.Range("C10").Value = 124
$ws.Range("D10").Value = 102
$ws.Range("E10").Value = 38
$ws.Range("B11").Value = 41
$ws.Range("C11").Value = 175
$ws.Range("D11").Value = 53
$ws.Range("E11").Value = 83

$ws = $wb.Worksheets.Item("ForcastedStandardDeviation")
$ws.Range("B2").Value = 1.8
$ws.Range("C2").Value = 3.524999999999999
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 2.774999999999999
$ws.Range("B3").Value = 3.562499999999999
$ws.Range("C3").Value = 10.2125
$ws.Range("D3").Value = 3.277499999999999
$ws.Range("E3").Value = 1.235
$ws.Range("B4").Value = 0.8807499999999997
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 5.690999999999998
$ws.Range("E4").Value = 0
$ws.Range("B5").Value = 7.135924999999999
$ws.Range("D5").Value = 7.651775
$ws.Range("E5").Value = 3.26705
$ws.Range("B6").Value = 5.323629999999999
$ws.Range("C6").Value = 3.583212499999999
$ws.Range("D6").Value = 8.394954999999998
$ws.Range("E6").Value = 8.190199999999999
$ws.Range("B7").Value = 8.668341499999999
$ws.Range("C7").Value = 5.271288749999999
$ws.Range("D7").Value = 2.342795
$ws.Range("E7").Value = 22.13941275
$ws.Range("B8").Value = 5.999585649999998
$ws.Range("C8").Value = 14.477261025
$ws.Range("D8").Value = 5.347456774999999
$ws.Range("E8").Value = 6.130011424999998
$ws.Range("B9").Value = 6.834393479999999
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 17.7978996875
$ws.Range("E9").Value = 3.132430345
$ws.Range("B10").Value = 12.864169731
$ws.Range("C10").Value = 18.989964841
$ws.Range("D10").Value = 15.6207775305
$ws.Range("E10").Value = 5.819505354499999
$ws.Range("B11").Value = 6.676045988975
$ws.Range("C11").Value = 28.495318245625
$ws.Range("D11").Value = 8.630010668675
$ws.Range("E11").Value = 13.514922367925

$ws = $wb.Worksheets.Item("Capacity")
$ws.Range("B2").Value = 1176
$ws.Range("B3").Value = 3384
$ws.Range("B4").Value = 1330
$ws.Range("B5").Value = 1902
$ws.Range("B6").Value = 4230
$ws.Range("B7").Value = 7548
$ws.Range("B8").Value = 3384
$ws.Range("B9").Value = 2598
$ws.Range("B10").Value = 1692
$ws.Range("B11").Value = 8198.999999999998
$ws.Range("B12").Value = 8198.999999999998
$ws.Range("B13").Value = 7548
$ws.Range("B14").Value = 4290
$ws.Range("B15").Value = 846
$ws.Range("B16").Value = 846
$ws.Range("B17").Value = 5196
$ws.Range("B18").Value = 9435

$ws = $wb.Worksheets.Item("ProcessingTime")
$ws.Range("B2").Value = 2
$ws.Range("C3").Value = 4
$ws.Range("D4").Value = 2
$ws.Range("F6").Value = 5
$ws.Range("H8").Value = 4
$ws.Range("I9").Value = 2
$ws.Range("K11").Value = 3
$ws.Range("L12").Value = 3
$ws.Range("M13").Value = 4
$ws.Range("N14").Value = 2
$ws.Range("O15").Value = 1
$ws.Range("P16").Value = 1
$ws.Range("Q17").Value = 4
$ws.Range("R18").Value = 5
